$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.568.73'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.760.65'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.00'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.11'
$ws.Range("E6").Value = '  +2.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.760.02'
$ws.Range("E7").Value = '  -1.78%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("E13").Value = '  +4.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.81'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.393.61'
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.759.88'
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.92'
$ws.Range("E17").Value = '  +4.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.655.42'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.56'
$ws.Range("E21").Value = '  -3.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '469.37'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.722'
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("E24").Value = '  -6.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.84'
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.23'
$ws.Range("E27").Value = '  +1.63%  '
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.907.48'
$ws.Range("E31").Value = '  -1.66%  '
$ws.Range("E32").Value = '  +1.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.26'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.43'
$ws.Range("E34").Value = '  -1.79%  '
$ws.Range("E35").Value = '  -4.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.727.91'
$ws.Range("E36").Value = '  -1.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.86'
$ws.Range("E37").Value = '  +7.85%  '
$ws.Range("E38").Value = '  +1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.91'
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("E40").Value = '  -1.37%  '
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.96'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.82'
$ws.Range("E47").Value = '  -2.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '398.24'
$ws.Range("E48").Value = '  -4.81%  '
$ws.Range("E49").Value = '  -8.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.51'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0355'
$ws.Range("E51").Value = '  -0.09%  '
